$d = $word.ActiveDocument

# 1) Remove the whole "1. Termo de Responsabilidade..." paragraph, the blank
#    paragraph right after it, and the "2. A usuária deverá fazer cadastro..."
#    paragraph (three consecutive paragraphs, including their paragraph marks).
$startFind = $d.Content.Find
$ok1 = $startFind.Execute("1. Termo de Responsabilidade", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPara = $d.Range($startFind.Parent.Start, $startFind.Parent.Start).Paragraphs(1)

$endFind = $d.Content.Find
$ok2 = $endFind.Execute("2. A usuária deverá fazer cadastro na Central de Segurança (PIÁ).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPara = $d.Range($endFind.Parent.Start, $endFind.Parent.Start).Paragraphs(1)

$delRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$delRange.Delete()

# 2) Drop the leading "3. " numbering from the "Comprovação de vínculo..."
#    paragraph, leaving the rest of the sentence/formatting untouched.
$ok3 = $d.Content.Find.Execute("3. Comprovação de vínculo profissional", $true, $false, $false, $false, $false, $true, 1, $false, "Comprovação de vínculo profissional", 2)
